$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @(
    "Liêm Trinh tọa thủ cung Mệnh ở Tý gặp Kình Dương, Đà La, Địa Không, Địa Kiếp, Hoả Tinh, Linh Tinh, Hoá Kỵ, Thiên Hình",
    "Liêm Trinh tọa thủ cung Mệnh ở Ngọ gặp Kình Dương, Đà La, Địa Không, Địa Kiếp, Hoả Tinh, Linh Tinh, Hoá Kỵ, Thiên Hình",
    "Liêm Trinh tọa thủ cung Mệnh ở Dần gặp Kình Dương, Đà La, Địa Không, Địa Kiếp, Hoả Tinh, Linh Tinh, Hoá Kỵ, Thiên Hình",
    "Liêm Trinh tọa thủ cung Mệnh ở Thân gặp Kình Dương, Đà La, Địa Không, Địa Kiếp, Hoả Tinh, Linh Tinh, Hoá Kỵ, Thiên Hình",
    "Liêm Trinh tọa thủ cung Mệnh ở Thìn gặp Kình Dương, Đà La, Địa Không, Địa Kiếp, Hoả Tinh, Linh Tinh, Hoá Kỵ, Thiên Hình",
    "Liêm Trinh tọa thủ cung Mệnh ở Tuất gặp Kình Dương, Đà La, Địa Không, Địa Kiếp, Hoả Tinh, Linh Tinh, Hoá Kỵ, Thiên Hình",
    "Liêm Trinh tọa thủ cung Mệnh ở Sửu gặp Kình Dương, Đà La, Địa Không, Địa Kiếp, Hoả Tinh, Linh Tinh, Hoá Kỵ, Thiên Hình",
    "Liêm Trinh tọa thủ cung Mệnh ở Mùi gặp Kình Dương, Đà La, Địa Không, Địa Kiếp, Hoả Tinh, Linh Tinh, Hoá Kỵ, Thiên Hình"
)

$startRow = 122
for ($i = 0; $i -lt $newValues.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $newValues[$i]
    $ws.Cells.Item($row, 2).Value = $newValues[$i]
    $ws.Range($ws.Cells.Item($row, 1), $ws.Cells.Item($row, 2)).Interior.Color = 65535
}

$lastRow = $startRow + $newValues.Length - 1
$ws.Cells.Item($lastRow + 2, 2).Select()
